# Fill in the Introduction section with the report's opening paragraph
# plus two "Part" stub paragraphs, replacing the single empty paragraph
# that previously followed the "Introduction" heading.

$d = $word.ActiveDocument

# The document currently ends with an empty paragraph right after the
# "Introduction" heading — that is the paragraph we turn into the first
# new line of body text.
$p = $d.Paragraphs.Last
$p.Range.Text = "This report investigates two linear algebra-based approaches to modelling and interpreting high-dimensional image data. The report is divided into two sections corresponding to these tasks, outlining the methods used, results, and relevant visualisations."

# Add "Part I" as its own paragraph.
$p.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Last
$p2.Range.Text = "Part I"

# Add "Part II" as its own paragraph.
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last
$p3.Range.Text = "Part II"

Write-Host "Introduction body text inserted."
